$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# row => (new Price, new Volume)
$updates = @(
    @{ Row = 2; D = '65.335.83'; E = '  -1.57%  ' }
    @{ Row = 3; D = '3.424.08'; E = '  -4.73%  ' }
    @{ Row = 4; D = $null; E = '  +0.07%  ' }
    @{ Row = 5; D = '594.21'; E = '  -2.08%  ' }
    @{ Row = 6; D = '134.72'; E = '  -9.21%  ' }
    @{ Row = 7; D = '3.422.61'; E = '  -4.80%  ' }
    @{ Row = 8; D = $null; E = '  -0.09%  ' }
    @{ Row = 9; D = '0.489'; E = '  +0.20%  ' }
    @{ Row = 10; D = '7.41'; E = '  -6.13%  ' }
    @{ Row = 11; D = $null; E = '  -11.09%  ' }
    @{ Row = 12; D = '0.375'; E = '  -9.61%  ' }
    @{ Row = 13; D = '4.000.29'; E = '  -4.75%  ' }
    @{ Row = 14; D = '0.0000179'; E = '  -13.03%  ' }
    @{ Row = 15; D = '26.32'; E = '  -11.17%  ' }
    @{ Row = 16; D = '3.430.84'; E = '  -4.92%  ' }
    @{ Row = 17; D = '65.271.59'; E = '  -1.73%  ' }
    @{ Row = 18; D = $null; E = '  -3.06%  ' }
    @{ Row = 19; D = '9.95'; E = '  -10.55%  ' }
    @{ Row = 20; D = '5.71'; E = '  -9.81%  ' }
    @{ Row = 21; D = '13.61'; E = '  -8.81%  ' }
    @{ Row = 22; D = '389.56'; E = '  -8.07%  ' }
    @{ Row = 23; D = '73.04'; E = '  -6.97%  ' }
    @{ Row = 24; D = '0.541'; E = '  -11.60%  ' }
    @{ Row = 25; D = $null; E = '  -0.03%  ' }
    @{ Row = 26; D = '3.566.08'; E = '  -4.33%  ' }
    @{ Row = 27; D = $null; E = '  -12.93%  ' }
    @{ Row = 28; D = '1.00'; E = '  +0.16%  ' }
    @{ Row = 29; D = '2.24'; E = '  -9.90%  ' }
    @{ Row = 30; D = '7.13'; E = '  -13.95%  ' }
    @{ Row = 31; D = '8.13'; E = '  -13.33%  ' }
    @{ Row = 32; D = '3.426.75'; E = '  -4.53%  ' }
    @{ Row = 33; D = $null; E = '  -0.03%  ' }
    @{ Row = 34; D = '0.143'; E = '  -8.71%  ' }
    @{ Row = 35; D = '22.52'; E = '  -10.33%  ' }
    @{ Row = 36; D = '172.74'; E = '  -1.27%  ' }
    @{ Row = 37; D = $null; E = '  -14.20%  ' }
    @{ Row = 38; D = '6.80'; E = '  -12.48%  ' }
    @{ Row = 39; D = $null; E = '  -9.34%  ' }
    @{ Row = 40; D = '4.80'; E = '  -14.41%  ' }
    @{ Row = 41; D = '0.0765'; E = '  -10.19%  ' }
    @{ Row = 42; D = '0.810'; E = '  -8.55%  ' }
    @{ Row = 43; D = '43.59'; E = '  -5.00%  ' }
    @{ Row = 44; D = '1.00'; E = '  +0.28%  ' }
    @{ Row = 45; D = '4.37'; E = '  -16.09%  ' }
    @{ Row = 46; D = '1.61'; E = '  -13.12%  ' }
    @{ Row = 47; D = $null; E = '  -3.64%  ' }
    @{ Row = 48; D = '21.71'; E = '  -7.80%  ' }
    @{ Row = 49; D = '6.51'; E = '  -8.95%  ' }
    @{ Row = 50; D = '2.12'; E = '  -16.03%  ' }
    @{ Row = 51; D = '2.195.69'; E = '  -8.13%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # These columns hold plain text (e.g. "1.00", "0.810"); force text
        # format first so Excel's smart-entry parser doesn't coerce
        # numeric-looking strings into numbers (dropping trailing zeros, etc.)
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}
